# Apply November-data update to pycaret_tables.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet: compare_models ---
$ws = $wb.Worksheets.Item("compare_models")
$ws.Range("C2").Value = 4.9968
$ws.Range("D2").Value = 59.8245
$ws.Range("E2").Value = 7.6199
$ws.Range("F2").Value = 0.9117
$ws.Range("G2").Value = 0.1297
$ws.Range("H2").Value = 0.0936
$ws.Range("I2").Value = 0.082
$ws.Range("C3").Value = 6.118
$ws.Range("D3").Value = 82.0008
$ws.Range("E3").Value = 8.9213
$ws.Range("F3").Value = 0.8837
$ws.Range("G3").Value = 0.1429
$ws.Range("H3").Value = 0.1087
$ws.Range("I3").Value = 0.056
$ws.Range("C4").Value = 6.876
$ws.Range("D4").Value = 99.9609
$ws.Range("E4").Value = 9.9526
$ws.Range("F4").Value = 0.8537
$ws.Range("G4").Value = 0.168
$ws.Range("H4").Value = 0.1286
$ws.Range("I4").Value = 0.034
$ws.Range("C5").Value = 6.4988
$ws.Range("D5").Value = 106.8479
$ws.Range("E5").Value = 10.1882
$ws.Range("F5").Value = 0.8481
$ws.Range("G5").Value = 0.1549
$ws.Range("H5").Value = 0.1129
$ws.Range("I5").Value = 0.102
$ws.Range("A6").Value = "lr"
$ws.Range("B6").Value = "Linear Regression"
$ws.Range("C6").Value = 8.2555
$ws.Range("D6").Value = 116.1555
$ws.Range("E6").Value = 10.7242
$ws.Range("F6").Value = 0.8254
$ws.Range("G6").Value = 0.2079
$ws.Range("H6").Value = 0.1598
$ws.Range("I6").Value = 1.586
$ws.Range("A7").Value = "ada"
$ws.Range("B7").Value = "AdaBoost Regressor"
$ws.Range("C7").Value = 7.579
$ws.Range("D7").Value = 119.9101
$ws.Range("E7").Value = 10.8098
$ws.Range("F7").Value = 0.8293
$ws.Range("G7").Value = 0.1771
$ws.Range("H7").Value = 0.1397
$ws.Range("I7").Value = 0.052
$ws.Range("A8").Value = "ridge"
$ws.Range("B8").Value = "Ridge Regression"
$ws.Range("C8").Value = 8.2894
$ws.Range("D8").Value = 124.7097
$ws.Range("E8").Value = 10.9927
$ws.Range("F8").Value = 0.8185
$ws.Range("G8").Value = 0.1937
$ws.Range("H8").Value = 0.1557
$ws.Range("I8").Value = 0.02
$ws.Range("A9").Value = "br"
$ws.Range("B9").Value = "Bayesian Ridge"
$ws.Range("C9").Value = 8.3415
$ws.Range("D9").Value = 127.7446
$ws.Range("E9").Value = 11.1211
$ws.Range("F9").Value = 0.815
$ws.Range("G9").Value = 0.194
$ws.Range("H9").Value = 0.1565
$ws.Range("I9").Value = 0.02
$ws.Range("A10").Value = "knn"
$ws.Range("B10").Value = "K Neighbors Regressor"
$ws.Range("C10").Value = 7.4425
$ws.Range("D10").Value = 132.1817
$ws.Range("E10").Value = 11.4265
$ws.Range("F10").Value = 0.7977
$ws.Range("G10").Value = 0.1845
$ws.Range("H10").Value = 0.1376
$ws.Range("I10").Value = 0.022
$ws.Range("A11").Value = "huber"
$ws.Range("B11").Value = "Huber Regressor"
$ws.Range("C11").Value = 8.849
$ws.Range("D11").Value = 139.4408
$ws.Range("E11").Value = 11.6839
$ws.Range("F11").Value = 0.7879
$ws.Range("G11").Value = 0.2144
$ws.Range("H11").Value = 0.1697
$ws.Range("I11").Value = 0.028
$ws.Range("A12").Value = "dt"
$ws.Range("B12").Value = "Decision Tree Regressor"
$ws.Range("C12").Value = 7.5852
$ws.Range("D12").Value = 175.0359
$ws.Range("E12").Value = 12.7388
$ws.Range("F12").Value = 0.7221
$ws.Range("G12").Value = 0.1996
$ws.Range("H12").Value = 0.1248
$ws.Range("I12").Value = 0.02
$ws.Range("C13").Value = 9.2884
$ws.Range("D13").Value = 184.5868
$ws.Range("E13").Value = 13.4314
$ws.Range("F13").Value = 0.7385
$ws.Range("G13").Value = 0.215
$ws.Range("H13").Value = 0.1723
$ws.Range("I13").Value = 0.694
$ws.Range("C14").Value = 9.2974
$ws.Range("D14").Value = 197.8792
$ws.Range("E14").Value = 13.9155
$ws.Range("F14").Value = 0.7216
$ws.Range("G14").Value = 0.2235
$ws.Range("H14").Value = 0.1773
$ws.Range("A15").Value = "par"
$ws.Range("B15").Value = "Passive Aggressive Regressor"
$ws.Range("C15").Value = 13.054
$ws.Range("D15").Value = 278.8507
$ws.Range("E15").Value = 16.5023
$ws.Range("F15").Value = 0.5725
$ws.Range("G15").Value = 0.3811
$ws.Range("H15").Value = 0.2566
$ws.Range("I15").Value = 0.02
$ws.Range("A16").Value = "omp"
$ws.Range("B16").Value = "Orthogonal Matching Pursuit"
$ws.Range("C16").Value = 10.8279
$ws.Range("D16").Value = 275.6453
$ws.Range("E16").Value = 16.5513
$ws.Range("F16").Value = 0.5948
$ws.Range("G16").Value = 0.2679
$ws.Range("H16").Value = 0.2152
$ws.Range("I16").Value = 0.018
$ws.Range("C17").Value = 16.0314
$ws.Range("D17").Value = 384.9043
$ws.Range("E17").Value = 19.2115
$ws.Range("F17").Value = 0.4872
$ws.Range("G17").Value = 0.3396
$ws.Range("H17").Value = 0.3289
$ws.Range("I17").Value = 0.02
$ws.Range("C18").Value = 23.0206
$ws.Range("D18").Value = 751.0234
$ws.Range("E18").Value = 27.055
$ws.Range("F18").Value = -0.0271
$ws.Range("G18").Value = 0.4748
$ws.Range("H18").Value = 0.4905
$ws.Range("I18").Value = 0.02
$ws.Range("C19").Value = 1534.4668
$ws.Range("D19").Value = 12717300.2025
$ws.Range("E19").Value = 1853.9792
$ws.Range("F19").Value = -11817.093
$ws.Range("G19").Value = 2.1585
$ws.Range("H19").Value = 33.9087
$ws.Range("I19").Value = 0.022

# --- Sheet: tuned_1 ---
$ws = $wb.Worksheets.Item("tuned_1")
$ws.Range("B2").Value = 4.3932
$ws.Range("C2").Value = 52.4281
$ws.Range("D2").Value = 7.2407
$ws.Range("E2").Value = 0.9099
$ws.Range("F2").Value = 0.1294
$ws.Range("G2").Value = 0.0863
$ws.Range("B3").Value = 5.1252
$ws.Range("C3").Value = 62.7556
$ws.Range("D3").Value = 7.9218
$ws.Range("E3").Value = 0.9256
$ws.Range("F3").Value = 0.1255
$ws.Range("G3").Value = 0.0935
$ws.Range("B4").Value = 6.5977
$ws.Range("C4").Value = 112.8718
$ws.Range("D4").Value = 10.6241
$ws.Range("E4").Value = 0.8362
$ws.Range("F4").Value = 0.2168
$ws.Range("G4").Value = 0.1504
$ws.Range("B5").Value = 6.2897
$ws.Range("C5").Value = 90.2957
$ws.Range("D5").Value = 9.5024
$ws.Range("E5").Value = 0.9165
$ws.Range("F5").Value = 0.1081
$ws.Range("G5").Value = 0.087
$ws.Range("B6").Value = 4.7812
$ws.Range("C6").Value = 47.6051
$ws.Range("D6").Value = 6.8996
$ws.Range("E6").Value = 0.8922
$ws.Range("F6").Value = 0.1177
$ws.Range("G6").Value = 0.0861
$ws.Range("B7").Value = 5.4374
$ws.Range("C7").Value = 73.1913
$ws.Range("D7").Value = 8.4377
$ws.Range("E7").Value = 0.8961
$ws.Range("F7").Value = 0.1395
$ws.Range("G7").Value = 0.1006
$ws.Range("B8").Value = 0.8592
$ws.Range("C8").Value = 24.7458
$ws.Range("D8").Value = 1.4127
$ws.Range("E8").Value = 0.0319
$ws.Range("F8").Value = 0.0393
$ws.Range("G8").Value = 0.025

# --- Sheet: tuned_2 ---
$ws = $wb.Worksheets.Item("tuned_2")
$ws.Range("B2").Value = 4.9566
$ws.Range("C2").Value = 43.9117
$ws.Range("D2").Value = 6.6266
$ws.Range("E2").Value = 0.9246
$ws.Range("F2").Value = 0.1269
$ws.Range("G2").Value = 0.0976
$ws.Range("B3").Value = 5.1335
$ws.Range("C3").Value = 55.8462
$ws.Range("D3").Value = 7.473
$ws.Range("E3").Value = 0.9338
$ws.Range("F3").Value = 0.1184
$ws.Range("G3").Value = 0.0892
$ws.Range("B4").Value = 7.0881
$ws.Range("C4").Value = 101.4404
$ws.Range("D4").Value = 10.0718
$ws.Range("E4").Value = 0.8528
$ws.Range("F4").Value = 0.2202
$ws.Range("G4").Value = 0.1593
$ws.Range("B5").Value = 5.8984
$ws.Range("C5").Value = 63.3446
$ws.Range("D5").Value = 7.9589
$ws.Range("E5").Value = 0.9414
$ws.Range("F5").Value = 0.1018
$ws.Range("G5").Value = 0.0856
$ws.Range("B6").Value = 4.9023
$ws.Range("C6").Value = 51.0342
$ws.Range("D6").Value = 7.1438
$ws.Range("E6").Value = 0.8844
$ws.Range("F6").Value = 0.1328
$ws.Range("G6").Value = 0.0936
$ws.Range("B7").Value = 5.5958
$ws.Range("C7").Value = 63.1154
$ws.Range("D7").Value = 7.8548
$ws.Range("E7").Value = 0.9074
$ws.Range("F7").Value = 0.14
$ws.Range("G7").Value = 0.1051
$ws.Range("B8").Value = 0.8273
$ws.Range("C8").Value = 20.1813
$ws.Range("D8").Value = 1.1904
$ws.Range("E8").Value = 0.0337
$ws.Range("F8").Value = 0.0414
$ws.Range("G8").Value = 0.0274

# --- Sheet: tuned_3 ---
$ws = $wb.Worksheets.Item("tuned_3")
$ws.Range("B2").Value = 4.2929
$ws.Range("C2").Value = 34.1549
$ws.Range("D2").Value = 5.8442
$ws.Range("E2").Value = 0.9413
$ws.Range("F2").Value = 0.1136
$ws.Range("G2").Value = 0.0881
$ws.Range("B3").Value = 6.035
$ws.Range("C3").Value = 84.5171
$ws.Range("D3").Value = 9.1933
$ws.Range("E3").Value = 0.8999
$ws.Range("F3").Value = 0.1418
$ws.Range("G3").Value = 0.1088
$ws.Range("B4").Value = 5.7199
$ws.Range("C4").Value = 61.085
$ws.Range("D4").Value = 7.8157
$ws.Range("E4").Value = 0.9113
$ws.Range("F4").Value = 0.1666
$ws.Range("G4").Value = 0.1214
$ws.Range("B5").Value = 7.2851
$ws.Range("C5").Value = 105.452
$ws.Range("D5").Value = 10.269
$ws.Range("E5").Value = 0.9025
$ws.Range("F5").Value = 0.139
$ws.Range("G5").Value = 0.1095
$ws.Range("B6").Value = 5.8976
$ws.Range("C6").Value = 71.2386
$ws.Range("D6").Value = 8.4403
$ws.Range("E6").Value = 0.8387
$ws.Range("F6").Value = 0.1462
$ws.Range("G6").Value = 0.1085
$ws.Range("B7").Value = 5.8461
$ws.Range("C7").Value = 71.2895
$ws.Range("D7").Value = 8.3125
$ws.Range("E7").Value = 0.8987
$ws.Range("F7").Value = 0.1414
$ws.Range("G7").Value = 0.1072
$ws.Range("B8").Value = 0.9526
$ws.Range("C8").Value = 23.7704
$ws.Range("D8").Value = 1.4805
$ws.Range("E8").Value = 0.0335
$ws.Range("F8").Value = 0.0169
$ws.Range("G8").Value = 0.0107

# --- Sheet: tuned_4 ---
$ws = $wb.Worksheets.Item("tuned_4")
$ws.Range("B2").Value = 4.8016
$ws.Range("C2").Value = 54.4375
$ws.Range("D2").Value = 7.3782
$ws.Range("E2").Value = 0.9065
$ws.Range("F2").Value = 0.1314
$ws.Range("G2").Value = 0.0936
$ws.Range("B3").Value = 5.1733
$ws.Range("C3").Value = 64.7339
$ws.Range("D3").Value = 8.0457
$ws.Range("E3").Value = 0.9233
$ws.Range("F3").Value = 0.1296
$ws.Range("G3").Value = 0.0959
$ws.Range("B4").Value = 6.5292
$ws.Range("C4").Value = 107.0817
$ws.Range("D4").Value = 10.348
$ws.Range("E4").Value = 0.8446
$ws.Range("F4").Value = 0.2145
$ws.Range("G4").Value = 0.1466
$ws.Range("B5").Value = 6.3799
$ws.Range("C5").Value = 90.5117
$ws.Range("D5").Value = 9.5138
$ws.Range("E5").Value = 0.9163
$ws.Range("F5").Value = 0.113
$ws.Range("G5").Value = 0.092
$ws.Range("B6").Value = 5.1982
$ws.Range("C6").Value = 82.5625
$ws.Range("D6").Value = 9.0864
$ws.Range("E6").Value = 0.813
$ws.Range("F6").Value = 0.1487
$ws.Range("G6").Value = 0.1003
$ws.Range("B7").Value = 5.6164
$ws.Range("C7").Value = 79.8655
$ws.Range("D7").Value = 8.8744
$ws.Range("E7").Value = 0.8807
$ws.Range("F7").Value = 0.1474
$ws.Range("G7").Value = 0.1057
$ws.Range("B8").Value = 0.7002
$ws.Range("C8").Value = 18.638
$ws.Range("D8").Value = 1.0536
$ws.Range("E8").Value = 0.0439
$ws.Range("F8").Value = 0.0354
$ws.Range("G8").Value = 0.0207

# --- Sheet: tuned_5 ---
$ws = $wb.Worksheets.Item("tuned_5")
$ws.Range("B2").Value = 7.2356
$ws.Range("C2").Value = 83.6585
$ws.Range("D2").Value = 9.1465
$ws.Range("E2").Value = 0.8563
$ws.Range("F2").Value = 0.2103
$ws.Range("G2").Value = 0.1467
$ws.Range("B3").Value = 7.396
$ws.Range("C3").Value = 96.8294
$ws.Range("D3").Value = 9.8402
$ws.Range("E3").Value = 0.8853
$ws.Range("F3").Value = 0.1864
$ws.Range("G3").Value = 0.1543
$ws.Range("B4").Value = 8.9513
$ws.Range("C4").Value = 145.9504
$ws.Range("D4").Value = 12.081
$ws.Range("E4").Value = 0.7881
$ws.Range("F4").Value = 0.2427
$ws.Range("G4").Value = 0.1941
$ws.Range("B5").Value = 9.1696
$ws.Range("C5").Value = 129.2794
$ws.Range("D5").Value = 11.3701
$ws.Range("E5").Value = 0.8805
$ws.Range("F5").Value = 0.1757
$ws.Range("G5").Value = 0.1474
$ws.Range("B6").Value = 8.5248
$ws.Range("C6").Value = 125.0595
$ws.Range("D6").Value = 11.183
$ws.Range("E6").Value = 0.7168
$ws.Range("F6").Value = 0.2243
$ws.Range("G6").Value = 0.1565
$ws.Range("B7").Value = 8.2555
$ws.Range("C7").Value = 116.1554
$ws.Range("D7").Value = 10.7242
$ws.Range("E7").Value = 0.8254
$ws.Range("F7").Value = 0.2079
$ws.Range("G7").Value = 0.1598
$ws.Range("B8").Value = 0.7964
$ws.Range("C8").Value = 22.6663
$ws.Range("D8").Value = 1.0714
$ws.Range("E8").Value = 0.0644
$ws.Range("F8").Value = 0.0244
$ws.Range("G8").Value = 0.0176

# --- Sheet: blend_model ---
$ws = $wb.Worksheets.Item("blend_model")
$ws.Range("B2").Value = 4.0504
$ws.Range("C2").Value = 34.571
$ws.Range("D2").Value = 5.8797
$ws.Range("E2").Value = 0.9406
$ws.Range("F2").Value = 0.1103
$ws.Range("G2").Value = 0.0769
$ws.Range("B3").Value = 4.8233
$ws.Range("C3").Value = 49.7678
$ws.Range("D3").Value = 7.0546
$ws.Range("E3").Value = 0.941
$ws.Range("F3").Value = 0.1168
$ws.Range("G3").Value = 0.0911
$ws.Range("B4").Value = 6.2201
$ws.Range("C4").Value = 87.0246
$ws.Range("D4").Value = 9.3287
$ws.Range("E4").Value = 0.8737
$ws.Range("F4").Value = 0.2037
$ws.Range("G4").Value = 0.1424
$ws.Range("B5").Value = 6.2923
$ws.Range("C5").Value = 74.2009
$ws.Range("D5").Value = 8.614
$ws.Range("E5").Value = 0.9314
$ws.Range("F5").Value = 0.11
$ws.Range("G5").Value = 0.0926
$ws.Range("B6").Value = 4.5813
$ws.Range("C6").Value = 55.602
$ws.Range("D6").Value = 7.4567
$ws.Range("E6").Value = 0.8741
$ws.Range("F6").Value = 0.1326
$ws.Range("G6").Value = 0.0885
$ws.Range("B7").Value = 5.1935
$ws.Range("C7").Value = 60.2333
$ws.Range("D7").Value = 7.6667
$ws.Range("E7").Value = 0.9122
$ws.Range("F7").Value = 0.1347
$ws.Range("G7").Value = 0.0983
$ws.Range("B8").Value = 0.9033
$ws.Range("C8").Value = 18.4521
$ws.Range("D8").Value = 1.206
$ws.Range("E8").Value = 0.0314
$ws.Range("F8").Value = 0.0354
$ws.Range("G8").Value = 0.0227

# --- Sheet: stack_model ---
$ws = $wb.Worksheets.Item("stack_model")
$ws.Range("B2").Value = 3.973
$ws.Range("C2").Value = 30.5412
$ws.Range("D2").Value = 5.5264
$ws.Range("E2").Value = 0.9475
$ws.Range("F2").Value = 0.1055
$ws.Range("G2").Value = 0.0771
$ws.Range("B3").Value = 3.8158
$ws.Range("C3").Value = 25.3678
$ws.Range("D3").Value = 5.0367
$ws.Range("E3").Value = 0.9699
$ws.Range("F3").Value = 0.0897
$ws.Range("G3").Value = 0.0731
$ws.Range("B4").Value = 6.0112
$ws.Range("C4").Value = 80.3829
$ws.Range("D4").Value = 8.9656
$ws.Range("E4").Value = 0.8833
$ws.Range("F4").Value = 0.1875
$ws.Range("G4").Value = 0.1291
$ws.Range("B5").Value = 5.5345
$ws.Range("C5").Value = 58.5784
$ws.Range("D5").Value = 7.6537
$ws.Range("E5").Value = 0.9458
$ws.Range("F5").Value = 0.096
$ws.Range("G5").Value = 0.0813
$ws.Range("B6").Value = 3.7151
$ws.Range("C6").Value = 27.9359
$ws.Range("D6").Value = 5.2854
$ws.Range("E6").Value = 0.9367
$ws.Range("F6").Value = 0.0999
$ws.Range("G6").Value = 0.0672
$ws.Range("B7").Value = 4.6099
$ws.Range("C7").Value = 44.5612
$ws.Range("D7").Value = 6.4936
$ws.Range("E7").Value = 0.9367
$ws.Range("F7").Value = 0.1157
$ws.Range("G7").Value = 0.0856
$ws.Range("B8").Value = 0.965
$ws.Range("C8").Value = 21.5454
$ws.Range("D8").Value = 1.5476
$ws.Range("E8").Value = 0.0288
$ws.Range("F8").Value = 0.0362
$ws.Range("G8").Value = 0.0223

# --- Sheet: pred_blend ---
$ws = $wb.Worksheets.Item("pred_blend")
$ws.Range("C2").Value = 4.991
$ws.Range("D2").Value = 45.6895
$ws.Range("E2").Value = 6.7594
$ws.Range("F2").Value = 0.9077
$ws.Range("G2").Value = 0.112
$ws.Range("H2").Value = 0.0853

# --- Sheet: pred_stack ---
$ws = $wb.Worksheets.Item("pred_stack")
$ws.Range("C2").Value = 4.616
$ws.Range("D2").Value = 35.8788
$ws.Range("E2").Value = 5.9899
$ws.Range("F2").Value = 0.9275
$ws.Range("G2").Value = 0.097
$ws.Range("H2").Value = 0.0755

# --- Sheet: pred_final ---
$ws = $wb.Worksheets.Item("pred_final")
$ws.Range("C2").Value = 1.2967
$ws.Range("D2").Value = 5.0315
$ws.Range("E2").Value = 2.2431
$ws.Range("F2").Value = 0.9925
$ws.Range("G2").Value = 0.0553
$ws.Range("H2").Value = 0.025
